$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "List" column (E2:E4) so every row shares the same list of names
# instead of the previous per-row numeric lists.
$ws.Range("E2:E4").Value = "ayush;deepanshu;aryan;sumit;"
